$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: updated estimate value
$ws.Range("B2").Value = -0.00000006988130962781725

# Row 3
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = -0.0007063166991173719
$ws.Range("F3").Value = 0.0005234280742637722
$ws.Range("G3").Value = 0.0003433014305587195
$ws.Range("H3").Value = 0.0009409900507647346
$ws.Range("I3").Value = 0.0000004770647812532526
$ws.Range("J3").Value = 0.0000004497263649919632
$ws.Range("K3").Value = 0.002097823931623929
$ws.Range("L3").Value = 0.000000005225441820439758
$ws.Range("M3").Value = 0.000000003934983426307987

# Row 4
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0.252525252525252
$ws.Range("D4").Value = 0.247474747474747
$ws.Range("E4").Value = -0.006323887507055811
$ws.Range("F4").Value = 0.01301996818610838
$ws.Range("G4").Value = 0.008470320293531447
$ws.Range("H4").Value = 0.01756716503656709
$ws.Range("I4").Value = 0.0001789228060869995
$ws.Range("J4").Value = 0.000164526812178966
$ws.Range("K4").Value = 0.04302632986889662
$ws.Range("L4").Value = 0.000001709752113614707
$ws.Range("M4").Value = 0.000001247776857399013

# Row 5
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.252525252525252
$ws.Range("D5").Value = 0.247474747474747
$ws.Range("E5").Value = -0.006323887507055811
$ws.Range("F5").Value = 0.01301996818610838
$ws.Range("G5").Value = 0.008470320293531447
$ws.Range("H5").Value = 0.01756716503656709
$ws.Range("I5").Value = 0.0001789228060869995
$ws.Range("J5").Value = 0.000164526812178966
$ws.Range("K5").Value = 0.04302632986889662
$ws.Range("L5").Value = 0.000001709752113614707
$ws.Range("M5").Value = 0.000001247776857399013

# Row 6
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0

# Row 7
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0.1571023579585537
$ws.Range("F7").Value = 0.07512575127821451
$ws.Range("G7").Value = 0.07183542636234121
$ws.Range("H7").Value = 0.03726666135395587
$ws.Range("I7").Value = 0.000008187448177777457
$ws.Range("J7").Value = 0.000007781342782182737
$ws.Range("K7").Value = 0.1509626341426867
$ws.Range("L7").Value = 0.0000000000000000000000000000001109335647967048
$ws.Range("M7").Value = 0.0000000000000000000000000000001111716196567835

# Row 8
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0

# Row 9
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0.7600000000000002
$ws.Range("D9").Value = 0.6840000000000003
$ws.Range("E9").Value = 0.006997874510246208
$ws.Range("F9").Value = 0.7400460062287872
$ws.Range("G9").Value = -0.06472226579351217
$ws.Range("H9").Value = 0.9755162690526452
$ws.Range("I9").Value = 1.701166841299321
$ws.Range("J9").Value = -0.1230740437405012
$ws.Range("K9").Value = 0.1443999999999991
$ws.Range("L9").Value = 0.0000000000000000000000000000001733336949948512
$ws.Range("M9").Value = 0.0000000000000000000000000000001737056557137243

# Row 10
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0.7600000000000002
$ws.Range("D10").Value = 0.6840000000000003
$ws.Range("E10").Value = 0.007000127162957136
$ws.Range("F10").Value = 0.7402230114369622
$ws.Range("G10").Value = -0.06473684335022185
$ws.Range("H10").Value = 0.07064339511985977
$ws.Range("I10").Value = 0.008921153296147446
$ws.Range("J10").Value = -0.0006454172420543787
$ws.Range("K10").Value = 0.1443999999999991
$ws.Range("L10").Value = 0.0000000000000000000000000000001733336949948512
$ws.Range("M10").Value = 0.0000000000000000000000000000001737056557137243

# Row 11
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0.8847368421052633
$ws.Range("D11").Value = 0.796263157894737
$ws.Range("E11").Value = 0.007552768781085297
$ws.Range("F11").Value = 0.861713907358402
$ws.Range("G11").Value = -0.07536193467570843
$ws.Range("H11").Value = 0.08223791357097242
$ws.Range("I11").Value = 0.01208988145582584
$ws.Range("J11").Value = -0.0008746647083570679
$ws.Range("K11").Value = 0.1680999999999999
$ws.Range("L11").Value = 0.00000000000000000000000000000001232595164407831
$ws.Range("M11").Value = 0.00000000000000000000000000000001235240218408706

# Row 12
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0.1494619134827799
$ws.Range("F12").Value = 0.06019618153569198
$ws.Range("G12").Value = 0.05702181709899352
$ws.Range("H12").Value = 0.02395582196303462
$ws.Range("I12").Value = 0.000003110343802961913
$ws.Range("J12").Value = 0.000002910382855321007
$ws.Range("K12").Value = 0.1091977290837823
$ws.Range("L12").Value = 0.000000000000000000000000000001294995294605977
$ws.Range("M12").Value = 0.00000000000000000000000000000126631423015305

# Row 13
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.1436233330666425
$ws.Range("F13").Value = 0.1005881602729378
$ws.Range("G13").Value = 0.09123126571988853
$ws.Range("H13").Value = 0.1004196530123883
$ws.Range("I13").Value = 0.00004662032970098705
$ws.Range("J13").Value = 0.00003940689327751749
$ws.Range("K13").Value = 0.166987643059551
$ws.Range("L13").Value = 0.0000000000000000000000000000000007703719777548943
$ws.Range("M13").Value = 0.0000000000000000000000000000000007720251365054413

# Row 14
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0.5025125628140698
$ws.Range("D14").Value = 0.497487437185929
$ws.Range("E14").Value = 0.1390245789549215
$ws.Range("F14").Value = 0.1193383307848196
$ws.Range("G14").Value = 0.1076045980978656
$ws.Range("H14").Value = 0.08877800613908313
$ws.Range("I14").Value = 0.00003613158088627324
$ws.Range("J14").Value = 0.00002994404435883512
$ws.Range("K14").Value = 0.1760537530080297
$ws.Range("L14").Value = 0.0000000000000000000000000000004437342591868191
$ws.Range("M14").Value = 0.0000000000000000000000000000004446864786271342

